$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# New column widths for columns Q (17) and R (18)
$ws.Columns.Item(17).ColumnWidth = 30.5
$ws.Columns.Item(18).ColumnWidth = 24.666666666666668

# Header row (row 2) - new field names
$ws.Range("Q2").Value = "testLong:long"
$ws.Range("R2").Value = "testDouble:double"
$ws.Range("S2").Value = "testDou0:double"
$ws.Range("T2").Value = "testDou1"

# Row 3 - long + double sample data
$ws.Range("Q3").Value = 84798398239797
$ws.Range("Q3").NumberFormat = "0"

$bigDouble = 2.8937489237893399 * [Math]::Pow(10, 22)
$ws.Range("R3").Value = $bigDouble
$ws.Range("S3").Value = $bigDouble
$ws.Range("T3").Value = $bigDouble

# Row 4
$ws.Range("Q4").Value = 5
$ws.Range("R4").Value = 84798398239797
$ws.Range("R4").NumberFormat = "0"
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 2

# Row 5
$ws.Range("Q5").Value = 4
$ws.Range("R5").Value = 5
$ws.Range("S5").Value = 3

# Row 6
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 2

# Update selection to match the new used range
$ws.Range("A1:T7").Select()
